# Add fixed O sample BEAST trees
# Updates the "O" serotype rows (4-7) with new tree-height / clock-rate
# values coming from re-run BEAST analyses on the fixed samples, adjusts
# the two data-column widths, and restores the last active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (O, sample 1) ---------------------------------------------
$ws.Range("B4").Value = 154
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2.427E-3"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").Value = "[102.2586, 217.8057]"
$ws.Range("E4").Value = "[1.7203E-3, 3.2179E-3]"

# --- Row 5 (O, sample 2) ---------------------------------------------
$ws.Range("B5").Value = 140
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2.646E-3"
$ws.Range("C5").ClearFormats()
$ws.Range("D5").Value = "[92.2312, 195.4317]"
$ws.Range("E5").Value = "[1.8441E-3, 3.4851E-3]"

# --- Row 6 (O, sample 3) ---------------------------------------------
$ws.Range("B6").Value = 103
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "3.131E-3"
$ws.Range("C6").ClearFormats()
$ws.Range("E6").Value = "[2.4521E-3, 3.8491E-3]"
$ws.Range("D6").Value = "[82.6899, 126.053]"

# --- Row 7 (O, sample 4) ---------------------------------------------
$ws.Range("B7").Value = 142
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "2.664E-3"
$ws.Range("C7").ClearFormats()
$ws.Range("D7").Value = "[96.3805, 197.7668]"
$ws.Range("E7").Value = "[1.8559E-3, 3.4921E-3]"

# --- Row 11 (SAT1) keeps the same displayed values --------------------
$ws.Range("D11").Value = "[108.1109, 618.647]"
$ws.Range("E11").Value = "[9.6834E-4, 2.6446E-3]"

# --- Column width tweaks ------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 17.333333333333332
$ws.Columns.Item(3).ColumnWidth = 24.166666666666668

# --- Restore the last selected cell -------------------------------------
$ws.Range("E8").Select() | Out-Null
